# Applies the latest crypto price/volume snapshot to Sheet1 (columns B:E, rows 2-51).
# Column A (rank index) and the header row are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cryptoData = @(
     @{ Row = 2; Coin = 'Bitcoin'; Link = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; Price = '29.236.42'; Volume = '  -0.27%  ' }
    ,@{ Row = 3; Coin = 'Ethereum'; Link = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; Price = '1.840.25'; Volume = '  -0.07%  ' }
    ,@{ Row = 4; Coin = 'TetherUSD'; Link = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; Price = '0.9994'; Volume = '  +0.01%  ' }
    ,@{ Row = 5; Coin = 'BNB'; Link = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; Price = '241.40'; Volume = '  -0.60%  ' }
    ,@{ Row = 6; Coin = 'XRP'; Link = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; Price = '0.6696'; Volume = '  -2.51%  ' }
    ,@{ Row = 7; Coin = 'USDC'; Link = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; Price = '1.000'; Volume = '  +0.03%  ' }
    ,@{ Row = 8; Coin = 'Dogecoin'; Link = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; Price = '0.07424'; Volume = '  -1.41%  ' }
    ,@{ Row = 9; Coin = 'Cardano'; Link = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; Price = '0.2961'; Volume = '  -2.29%  ' }
    ,@{ Row = 10; Coin = 'Solana'; Link = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; Price = '22.80'; Volume = '  -2.13%  ' }
    ,@{ Row = 11; Coin = 'TRON'; Link = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Price = '0.07719'; Volume = '  +0.85%  ' }
    ,@{ Row = 12; Coin = 'WrappedEther'; Link = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; Price = '1.920.21'; Volume = '  +4.36%  ' }
    ,@{ Row = 13; Coin = 'Polkadot'; Link = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Price = '5.030'; Volume = '  -1.10%  ' }
    ,@{ Row = 14; Coin = 'Polygon'; Link = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; Price = '0.6786'; Volume = '  -1.22%  ' }
    ,@{ Row = 15; Coin = 'Litecoin'; Link = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Price = '86.32'; Volume = '  -3.08%  ' }
    ,@{ Row = 16; Coin = 'Uniswap'; Link = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; Price = '6.196'; Volume = '  -1.67%  ' }
    ,@{ Row = 17; Coin = 'ShibaInu'; Link = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Price = '0.000008264'; Volume = '  +0.34%  ' }
    ,@{ Row = 18; Coin = 'WrappedBTC'; Link = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; Price = '28.657.55'; Volume = '  -2.23%  ' }
    ,@{ Row = 19; Coin = 'BitcoinCash'; Link = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; Price = '228.73'; Volume = '  -1.93%  ' }
    ,@{ Row = 20; Coin = 'Avalanche'; Link = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; Price = '12.55'; Volume = '  -0.42%  ' }
    ,@{ Row = 21; Coin = 'Dai'; Link = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; Price = '0.9996'; Volume = '  +0.03%  ' }
    ,@{ Row = 22; Coin = 'Chainlink'; Link = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Price = '7.191'; Volume = '  -4.04%  ' }
    ,@{ Row = 23; Coin = 'BinanceUSD'; Link = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; Price = '1.000'; Volume = '  +0.01%  ' }
    ,@{ Row = 24; Coin = 'Monero'; Link = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Price = '160.35'; Volume = '  +0.17%  ' }
    ,@{ Row = 25; Coin = 'Cosmos'; Link = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Price = '8.711'; Volume = '  -1.53%  ' }
    ,@{ Row = 26; Coin = 'Stellar'; Link = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Price = '0.1409'; Volume = '  -3.54%  ' }
    ,@{ Row = 27; Coin = 'EthereumClassic'; Link = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; Price = '18.03'; Volume = '  -0.33%  ' }
    ,@{ Row = 28; Coin = 'PancakeSwap'; Link = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Price = '1.508'; Volume = '  -0.54%  ' }
    ,@{ Row = 29; Coin = 'Filecoin'; Link = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Price = '4.200'; Volume = '  -0.75%  ' }
    ,@{ Row = 30; Coin = 'InternetComputer(DFINITY)'; Link = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Price = '4.091'; Volume = '  -1.25%  ' }
    ,@{ Row = 31; Coin = 'Toncoin'; Link = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Price = '1.187'; Volume = '  -1.55%  ' }
    ,@{ Row = 32; Coin = 'Hedera'; Link = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; Price = '0.05367'; Volume = '  +4.40%  ' }
    ,@{ Row = 33; Coin = 'LidoDAOToken'; Link = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; Price = '1.877'; Volume = '  +1.61%  ' }
    ,@{ Row = 34; Coin = 'ImmutableX'; Link = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Price = '0.7540'; Volume = '  -2.62%  ' }
    ,@{ Row = 35; Coin = 'ARBITRUM'; Link = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; Price = '1.140'; Volume = '  +0.12%  ' }
    ,@{ Row = 36; Coin = 'HuobiToken'; Link = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; Price = '2.676'; Volume = '  +0.19%  ' }
    ,@{ Row = 37; Coin = 'Maker'; Link = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; Price = '1.329.72'; Volume = '  +3.37%  ' }
    ,@{ Row = 38; Coin = 'VeChain'; Link = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Price = '0.01804'; Volume = '  -2.20%  ' }
    ,@{ Row = 39; Coin = 'MXToken'; Link = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Price = '2.732'; Volume = '  +1.24%  ' }
    ,@{ Row = 40; Coin = 'TrustWalletToken'; Link = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Price = '0.9209'; Volume = '  -2.37%  ' }
    ,@{ Row = 41; Coin = 'FraxShare'; Link = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Price = '5.977'; Volume = '  +5.55%  ' }
    ,@{ Row = 42; Coin = 'PaxDollar'; Link = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; Price = '0.9997'; Volume = '  -0.02%  ' }
    ,@{ Row = 43; Coin = 'Quant'; Link = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Price = '103.09'; Volume = '  -2.47%  ' }
    ,@{ Row = 44; Coin = 'XinFinNetwork'; Link = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'; Price = '0.07886'; Volume = '  +12.40%  ' }
    ,@{ Row = 45; Coin = 'BabyDogeCoin'; Link = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; Price = '0.00000000123'; Volume = '  +3.68%  ' }
    ,@{ Row = 46; Coin = 'Mantle'; Link = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; Price = '0.5164'; Volume = '  -0.66%  ' }
    ,@{ Row = 47; Coin = 'Aave'; Link = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; Price = '64.10'; Volume = '  +1.36%  ' }
    ,@{ Row = 48; Coin = 'RenderToken'; Link = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Price = '1.764'; Volume = '  -0.54%  ' }
    ,@{ Row = 49; Coin = 'EnergySwap'; Link = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Price = '9.237'; Volume = '  -4.62%  ' }
    ,@{ Row = 50; Coin = 'Cronos'; Link = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Price = '0.05923'; Volume = '  -0.05%  ' }
    ,@{ Row = 51; Coin = 'Aptos'; Link = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Price = '6.877'; Volume = '  -0.37%  ' }
)

foreach ($item in $cryptoData) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Coin
    $ws.Range("C$r").Value = $item.Link

    # Price column: force text so strings like "1.000" or "241.40" are not
    # silently coerced into numbers (which would drop the formatting/precision),
    # then restore the default (unstyled) cell style used by the rest of the sheet.
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $item.Price
    $ws.Range("D$r").Style = "Normal"

    $ws.Range("E$r").Value = $item.Volume
}
